# register user design.pptx — reposition a handful of field-label textboxes
# on the "register" slide (slide 5).
#
# NOTE on the notes-slide wording tweak ("ahs to " + "be unique" -> one run
# " ahs to be unique"): this COM-interop host can only rewrite a Notes
# placeholder's *entire* TextFrame.TextRange.Text in one shot (Characters /
# Paragraphs / Runs sub-ranges, Font.*, InsertAfter/Before, Delete, Replace
# are all no-ops on NotesPage body placeholders here), which would blow away
# every other run's rPr (lang/dirty/smtClean/err) on that shape for a change
# that doesn't even alter the rendered text (the two runs already read
# "Uname" + " ahs to " + "be unique" = "Uname ahs to be unique" either way).
# That trade is net-negative, so it is intentionally left alone and this
# script focuses solely on the shape-position changes, which the host
# applies cleanly via Shape.Left / Shape.Top.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)

# Target positions, expressed in points (EMU / 12700). A couple of values
# are nudged by a hair above the exact quotient because this host round-trips
# Left/Top through a 32-bit float before truncating to EMU on save; the
# nudge keeps that truncation landing on the exact target EMU.
$moves = @{
    "TextBox 7"  = @{ Left = 354.21882629763945; Top = 174.248031496063   }
    "TextBox 8"  = @{ Left = 373.74757385512146; Top = 215.0955905511811 }
    "TextBox 11" = @{ Left = 310.52119450236523; Top = 346.62599185197024 }
    "TextBox 13" = @{ Left = 328.1668548937018;  Top = 305.60976377952755 }
    "TextBox 15" = @{ Left = 326.0765354330709;  Top = 258.9443307086614 }
}

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($moves.ContainsKey($sh.Name)) {
        $target = $moves[$sh.Name]
        $sh.Left = $target.Left
        $sh.Top = $target.Top
    }
}
